$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 116: Steven / Wily Hyena workout on 2024-06-26 ---
$ws.Cells.Item(116, 1).Value = "Steven"
$ws.Cells.Item(116, 2).Value = 45469
$ws.Cells.Item(116, 3).Value = "Workout"
$ws.Cells.Item(116, 4).Value = 42
$ws.Cells.Item(116, 5).Value = 0
$ws.Cells.Item(116, 6).Value = 0
$ws.Cells.Item(116, 7).Value = 3
$ws.Cells.Item(116, 8).Value = 19
$ws.Cells.Item(116, 9).Value = 19
$ws.Cells.Item(116, 10).Value = 1
$ws.Cells.Item(116, 11).Value = 0
$ws.Cells.Item(116, 12).Value = "Wily Hyena"
$ws.Cells.Item(116, 13).Value = 3

# --- Row 117: Steven / Wily Hyena run on 2024-06-27 ---
$ws.Cells.Item(117, 1).Value = "Steven"
$ws.Cells.Item(117, 2).Value = 45470
$ws.Cells.Item(117, 3).Value = "Run"
$ws.Cells.Item(117, 4).Value = 31
$ws.Cells.Item(117, 5).Value = 3
$ws.Cells.Item(117, 6).Value = 92
$ws.Cells.Item(117, 7).Value = 2
$ws.Cells.Item(117, 8).Value = 19
$ws.Cells.Item(117, 9).Value = 10
$ws.Cells.Item(117, 10).Value = 0
$ws.Cells.Item(117, 11).Value = 0
$ws.Cells.Item(117, 12).Value = "Wily Hyena"
$ws.Cells.Item(117, 13).Value = 3

# Apply the same date format/style used by the rest of column B (reuse the
# existing style index instead of letting a brand-new numFmt get created).
$ws.Range("B115").Copy()
$ws.Range("B116:B117").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Move the active selection the way Excel would after typing the new rows:
# one row below the freshly-entered data (matches the saved view state).
$ws.Range("A118").Select()
